$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7944
$ws1.Range("F3").Value = 77
$ws1.Range("F5").Value = 62
$ws1.Range("F6").Value = 742
$ws1.Range("F7").Value = 1337
$ws1.Range("F8").Value = 219
$ws1.Range("F9").Value = 25
$ws1.Range("F10").Value = 185
$ws1.Range("F11").Value = 41

# Sheet "全部类型" (All types) - row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7944
$ws4.Range("F3").Value = 77
$ws4.Range("F5").Value = 62
$ws4.Range("F6").Value = 742
$ws4.Range("F7").Value = 1337
$ws4.Range("F8").Value = 219
$ws4.Range("F10").Value = 25
$ws4.Range("F11").Value = 185
$ws4.Range("F12").Value = 41
